$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q3" right before the current "2022-Q2"
#    sheet (sheet index 2), so sheet order becomes:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$q3.Name = "2022-Q3"

# Header row (row 1)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Column A index counters are numeric (0-based row counters)
$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1
$q3.Range("A4").Value = 2
$q3.Range("A5").Value = 3

# Text-like numeric columns B,C,D,E,F,G must stay TEXT (not auto-converted
# to numbers / losing leading zeros), so force the number format to Text
# before assigning the values.
$textCols = $q3.Range("B2:G5")
$textCols.NumberFormat = "@"

$q3.Range("B2").Value = "501054"
$q3.Range("C2").Value = "东方红睿泽三年定期开放灵活配置混合A"
$q3.Range("D2").Value = "100.44"
$q3.Range("E2").Value = "95.37"
$q3.Range("F2").Value = "3.09"
$q3.Range("G2").Value = "3.1036"
$q3.Range("H2").Value = 9

$q3.Range("B3").Value = "009576"
$q3.Range("C3").Value = "东方红智远三年持有期混合"
$q3.Range("D3").Value = "62.81"
$q3.Range("E3").Value = "92.83"
$q3.Range("F3").Value = "3.01"
$q3.Range("G3").Value = "1.8906"
$q3.Range("H3").Value = 8

$q3.Range("B4").Value = "169104"
$q3.Range("C4").Value = "东方红睿满沪港深灵活配置混合（LOF）"
$q3.Range("D4").Value = "42.07"
$q3.Range("E4").Value = "91.66"
$q3.Range("F4").Value = "3.87"
$q3.Range("G4").Value = "1.6281"
$q3.Range("H4").Value = 6

$q3.Range("B5").Value = "011032"
$q3.Range("C5").Value = "东方红睿泽三年定期开放灵活配置混合C"
$q3.Range("D5").Value = "0.32"
$q3.Range("E5").Value = "95.37"
$q3.Range("F5").Value = "3.09"
$q3.Range("G5").Value = "0.0099"
$q3.Range("H5").Value = 9

# Strip the lingering "@" text-format style off the cells we just forced to
# text (copy the plain/default format from an untouched cell over them),
# while keeping their stored string values & type intact.
$q3.Range("Z1").Copy()
$textCols.PasteSpecial(-4122)

# Apply the bold/bordered/centered header style to row 1 and to the column A
# index cells, matching the other quarter sheets - by copying format only
# from a cell (on an existing sheet) that already carries that style. Look
# the sheet back up by name (sheet object refs are position-bound in this
# engine, and the "2022-Q2" sheet shifted position when $q3 was inserted).
$q2 = $wb.Worksheets.Item("2022-Q2")
$styleSrcHeader = $q2.Range("B1")
$styleSrcIndex = $q2.Range("A2")

$styleSrcHeader.Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$styleSrcIndex.Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new data row for 2022-Q3
#    right after the header, shifting the existing rows down, and
#    renumber the index column A.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Read the existing 5 data rows (B2:D6) before we overwrite them, so we can
# shift them down by one row.
$existing = @()
for ($r = 2; $r -le 6; $r++) {
    $b = $total.Cells.Item($r, 2).Value2
    $c = $total.Cells.Item($r, 3).Value2
    $d = $total.Cells.Item($r, 4).Value2
    $existing += , @($b, $c, $d)
}

# Row 2 becomes the new 2022-Q3 entry.
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 6.63

# Rows 3-7 get the previous rows 2-6 shifted down by one.
for ($i = 0; $i -lt $existing.Length; $i++) {
    $r = $i + 3
    $total.Cells.Item($r, 2).Value = $existing[$i][0]
    $total.Cells.Item($r, 3).Value = $existing[$i][1]
    $total.Cells.Item($r, 4).Value = $existing[$i][2]
}

# Row 7 is brand new - give its index cell (A7) the same bold/bordered style
# used by the other column-A index cells (A2:A6).
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

# Renumber column A (0-based sequential counter) for rows 2-7.
for ($r = 2; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
